$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1595.7073
$ws_ALC.Range("J17").Value = 1595.7073
$ws_ALC.Range("L17").Value = 4787.1219
$ws_ALC.Range("N17").Value = -5123.1219
$ws_ALC.Range("H19").Value = 1337
$ws_ALC.Range("I19").Value = 1374.5
$ws_ALC.Range("J19").Value = 1299.5
$ws_ALC.Range("K19").Value = 1374.5
$ws_ALC.Range("L19").Value = 1299.5
$ws_ALC.Range("M19").Value = -1199.5
$ws_ALC.Range("N19").Value = -1649.5
$ws_ALC.Range("H64").Value = 6399.1333
$ws_ALC.Range("I64").Value = 3609
$ws_ALC.Range("J64").Value = 8259.223
$ws_ALC.Range("K64").Value = 3609
$ws_ALC.Range("L64").Value = 8259.223
$ws_ALC.Range("M64").Value = -3361
$ws_ALC.Range("N64").Value = -8755.223
$ws_ALC.Range("H67").Value = 6399.1333
$ws_ALC.Range("I67").Value = 3609
$ws_ALC.Range("J67").Value = 8259.223
$ws_ALC.Range("K67").Value = 3609
$ws_ALC.Range("L67").Value = 8259.223
$ws_ALC.Range("M67").Value = -2751
$ws_ALC.Range("N67").Value = -9975.223
$ws_ALC.Range("H74").Value = 13519
$ws_ALC.Range("I74").Value = 14570.9
$ws_ALC.Range("K74").Value = 14570.9
$ws_ALC.Range("M74").Value = -13634.9
$ws_ALC.Range("H76").Value = 3423.077
$ws_ALC.Range("I76").Value = 3095.238
$ws_ALC.Range("K76").Value = 3095.238
$ws_ALC.Range("M76").Value = -2780.238
$ws_ALC.Range("H77").Value = 13519
$ws_ALC.Range("I77").Value = 14570.9
$ws_ALC.Range("K77").Value = 72854.5
$ws_ALC.Range("M77").Value = -68174.5
$ws_ALC.Range("H79").Value = 3423.077
$ws_ALC.Range("I79").Value = 3095.238
$ws_ALC.Range("K79").Value = 3095.238
$ws_ALC.Range("M79").Value = -2003.238
$ws_ALC.Range("H96").Value = 667546.9399999999
$ws_ALC.Range("I96").Value = 1111715.6
$ws_ALC.Range("J96").Value = 1293.8334
$ws_ALC.Range("K96").Value = 3335146.8
$ws_ALC.Range("L96").Value = 3881.5002
$ws_ALC.Range("M96").Value = -3333773.8
$ws_ALC.Range("N96").Value = -6627.5002
$ws_ALC.Range("H98").Value = 8094.227
$ws_ALC.Range("I98").Value = 12176.786
$ws_ALC.Range("K98").Value = 12176.786
$ws_ALC.Range("M98").Value = -10678.786
$ws_ALC.Range("H122").Value = 8094.227
$ws_ALC.Range("I122").Value = 12176.786
$ws_ALC.Range("K122").Value = 36530.358
$ws_ALC.Range("M122").Value = -34080.358
$ws_ALC.Range("H137").Value = 3620.818
$ws_ALC.Range("J137").Value = 4999.5
$ws_ALC.Range("L137").Value = 14998.5
$ws_ALC.Range("N137").Value = -20098.5
$ws_ALC.Range("H138").Value = 3052.6
$ws_ALC.Range("I138").Value = 2726.2856
$ws_ALC.Range("K138").Value = 8178.8568
$ws_ALC.Range("M138").Value = -3038.8568
$ws_ALC.Range("H141").Value = 4038
$ws_ALC.Range("I141").Value = 3831.9285
$ws_ALC.Range("K141").Value = 11495.7855
$ws_ALC.Range("M141").Value = -6315.7855

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 7794.25
$ws_ARM.Range("I32").Value = 7194.35
$ws_ARM.Range("K32").Value = 7194.35
$ws_ARM.Range("M32").Value = -6907.35
$ws_ARM.Range("H88").Value = 3650.5
$ws_ARM.Range("I88").Value = 2953
$ws_ARM.Range("J88").Value = 3999.25
$ws_ARM.Range("K88").Value = 2953
$ws_ARM.Range("L88").Value = 3999.25
$ws_ARM.Range("M88").Value = -2547
$ws_ARM.Range("N88").Value = -4811.25
$ws_ARM.Range("H91").Value = 3650.5
$ws_ARM.Range("I91").Value = 2953
$ws_ARM.Range("J91").Value = 3999.25
$ws_ARM.Range("K91").Value = 2953
$ws_ARM.Range("L91").Value = 3999.25
$ws_ARM.Range("M91").Value = -1549
$ws_ARM.Range("N91").Value = -6807.25
$ws_ARM.Range("H135").Value = 90285.39999999999
$ws_ARM.Range("J135").Value = 90285.39999999999
$ws_ARM.Range("L135").Value = 90285.39999999999
$ws_ARM.Range("N135").Value = -100425.4

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H102").Value = 7894
$ws_BSM.Range("I102").Value = 7894
$ws_BSM.Range("K102").Value = 7894
$ws_BSM.Range("M102").Value = -4649

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H5").Value = 2271.4167
$ws_CRP.Range("I5").Value = 2618.625
$ws_CRP.Range("J5").Value = 1577
$ws_CRP.Range("K5").Value = 2618.625
$ws_CRP.Range("L5").Value = 1577
$ws_CRP.Range("M5").Value = -2506.625
$ws_CRP.Range("N5").Value = -1801
$ws_CRP.Range("H88").Value = 19889.75
$ws_CRP.Range("I88").Value = 20311
$ws_CRP.Range("J88").Value = 19749.334
$ws_CRP.Range("K88").Value = 20311
$ws_CRP.Range("L88").Value = 19749.334
$ws_CRP.Range("M88").Value = -19905
$ws_CRP.Range("N88").Value = -20561.334
$ws_CRP.Range("H91").Value = 19889.75
$ws_CRP.Range("I91").Value = 20311
$ws_CRP.Range("J91").Value = 19749.334
$ws_CRP.Range("K91").Value = 20311
$ws_CRP.Range("L91").Value = 19749.334
$ws_CRP.Range("M91").Value = -18907
$ws_CRP.Range("N91").Value = -22557.334
$ws_CRP.Range("H93").Value = 49900
$ws_CRP.Range("J93").Value = 0
$ws_CRP.Range("L93").Value = 0
$ws_CRP.Range("N93").ClearContents()
$ws_CRP.Range("H105").Value = 2321.2354
$ws_CRP.Range("I105").Value = 1746.5
$ws_CRP.Range("K105").Value = 1746.5
$ws_CRP.Range("M105").Value = 0.5

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H12").Value = 136.77777
$ws_CUL.Range("I12").Value = 32.714287
$ws_CUL.Range("K12").Value = 98.142861
$ws_CUL.Range("M12").Value = 74.857139
$ws_CUL.Range("H46").Value = 12813963
$ws_CUL.Range("I46").Value = 11616290
$ws_CUL.Range("J46").Value = 20000000
$ws_CUL.Range("K46").Value = 34848870
$ws_CUL.Range("L46").Value = 60000000
$ws_CUL.Range("M46").Value = -34848779
$ws_CUL.Range("N46").Value = -60000182
$ws_CUL.Range("H114").Value = 557.25
$ws_CUL.Range("J114").Value = 709.6667
$ws_CUL.Range("L114").Value = 2129.0001
$ws_CUL.Range("N114").Value = -8637.000100000001
$ws_CUL.Range("H120").Value = 36135.56
$ws_CUL.Range("I120").Value = 17365.572
$ws_CUL.Range("K120").Value = 52096.716
$ws_CUL.Range("M120").Value = -47258.716
$ws_CUL.Range("H134").Value = 1546.3572
$ws_CUL.Range("I134").Value = 1203.7693
$ws_CUL.Range("K134").Value = 3611.3079
$ws_CUL.Range("M134").Value = 1458.6921
$ws_CUL.Range("H138").Value = 4598.7646
$ws_CUL.Range("I138").Value = 2029.8334
$ws_CUL.Range("J138").Value = 6000
$ws_CUL.Range("K138").Value = 6089.5002
$ws_CUL.Range("L138").Value = 18000
$ws_CUL.Range("M138").Value = -949.5002000000004
$ws_CUL.Range("N138").Value = -28280
$ws_CUL.Range("H140").Value = 61263.59
$ws_CUL.Range("I140").Value = 72677.21000000001
$ws_CUL.Range("J140").Value = 8000
$ws_CUL.Range("K140").Value = 218031.63
$ws_CUL.Range("L140").Value = 24000
$ws_CUL.Range("M140").Value = -212851.63
$ws_CUL.Range("N140").Value = -34360

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 5963.4443
$ws_GSM.Range("I80").Value = 4334.2
$ws_GSM.Range("J80").Value = 8000
$ws_GSM.Range("K80").Value = 4334.2
$ws_GSM.Range("L80").Value = 8000
$ws_GSM.Range("M80").Value = -3336.2
$ws_GSM.Range("N80").Value = -9996
$ws_GSM.Range("H83").Value = 5963.4443
$ws_GSM.Range("I83").Value = 4334.2
$ws_GSM.Range("J83").Value = 8000
$ws_GSM.Range("K83").Value = 21671
$ws_GSM.Range("L83").Value = 40000
$ws_GSM.Range("M83").Value = -16679
$ws_GSM.Range("N83").Value = -49984
$ws_GSM.Range("H113").Value = 7210.5264
$ws_GSM.Range("I113").Value = 4100
$ws_GSM.Range("J113").Value = 10666.667
$ws_GSM.Range("K113").Value = 4100
$ws_GSM.Range("L113").Value = 10666.667
$ws_GSM.Range("M113").Value = -1930
$ws_GSM.Range("N113").Value = -15006.667
$ws_GSM.Range("H122").Value = 89999.5
$ws_GSM.Range("I122").Value = 89999.5
$ws_GSM.Range("K122").Value = 269998.5
$ws_GSM.Range("M122").Value = -267548.5

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H10").Value = 5309
$ws_LTW.Range("J10").Value = 6561.25
$ws_LTW.Range("L10").Value = 6561.25
$ws_LTW.Range("N10").Value = -6841.25
$ws_LTW.Range("H16").Value = 463
$ws_LTW.Range("I16").Value = 467.66666
$ws_LTW.Range("J16").Value = 449
$ws_LTW.Range("K16").Value = 467.66666
$ws_LTW.Range("L16").Value = 449
$ws_LTW.Range("M16").Value = -297.66666
$ws_LTW.Range("N16").Value = -789
$ws_LTW.Range("H53").Value = 41173
$ws_LTW.Range("J53").Value = 50000
$ws_LTW.Range("L53").Value = 50000
$ws_LTW.Range("N53").Value = -51036
$ws_LTW.Range("H55").Value = 1474.762
$ws_LTW.Range("I55").Value = 243.81818
$ws_LTW.Range("J55").Value = 2828.8
$ws_LTW.Range("K55").Value = 243.81818
$ws_LTW.Range("L55").Value = 2828.8
$ws_LTW.Range("M55").Value = -70.81818000000001
$ws_LTW.Range("N55").Value = -3174.8
$ws_LTW.Range("H68").Value = 5991.88
$ws_LTW.Range("I68").Value = 4176.6924
$ws_LTW.Range("J68").Value = 7958.3335
$ws_LTW.Range("K68").Value = 4176.6924
$ws_LTW.Range("L68").Value = 7958.3335
$ws_LTW.Range("M68").Value = -3427.6924
$ws_LTW.Range("N68").Value = -9456.333500000001
$ws_LTW.Range("H71").Value = 5991.88
$ws_LTW.Range("I71").Value = 4176.6924
$ws_LTW.Range("J71").Value = 7958.3335
$ws_LTW.Range("K71").Value = 20883.462
$ws_LTW.Range("L71").Value = 39791.6675
$ws_LTW.Range("M71").Value = -17139.462
$ws_LTW.Range("N71").Value = -47279.6675
$ws_LTW.Range("H82").Value = 2767.1177
$ws_LTW.Range("J82").Value = 4199
$ws_LTW.Range("L82").Value = 4199
$ws_LTW.Range("N82").Value = -4921
$ws_LTW.Range("H85").Value = 2767.1177
$ws_LTW.Range("J85").Value = 4199
$ws_LTW.Range("L85").Value = 4199
$ws_LTW.Range("N85").Value = -6695
$ws_LTW.Range("H93").Value = 3859.2
$ws_LTW.Range("I93").Value = 361
$ws_LTW.Range("K93").Value = 361
$ws_LTW.Range("M93").Value = 887
$ws_LTW.Range("H132").Value = 2719.5122
$ws_LTW.Range("I132").Value = 2628.282
$ws_LTW.Range("J132").Value = 4498.5
$ws_LTW.Range("K132").Value = 7884.846
$ws_LTW.Range("L132").Value = 13495.5
$ws_LTW.Range("M132").Value = -5354.846
$ws_LTW.Range("N132").Value = -18555.5

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 2432.1428
$ws_WVR.Range("I107").Value = 2037.375
$ws_WVR.Range("K107").Value = 6112.125
$ws_WVR.Range("M107").Value = -4192.125
$ws_WVR.Range("H122").Value = 4590.476
$ws_WVR.Range("I122").Value = 2430.6155
$ws_WVR.Range("J122").Value = 8100.25
$ws_WVR.Range("K122").Value = 7291.8465
$ws_WVR.Range("L122").Value = 24300.75
$ws_WVR.Range("M122").Value = -4841.8465
$ws_WVR.Range("N122").Value = -29200.75
$ws_WVR.Range("H123").Value = 34997.5
$ws_WVR.Range("J123").Value = 34997.5
$ws_WVR.Range("L123").Value = 34997.5
$ws_WVR.Range("N123").Value = -44797.5
$ws_WVR.Range("H126").Value = 1994.5
$ws_WVR.Range("J126").Value = 0
$ws_WVR.Range("L126").Value = 0
$ws_WVR.Range("N126").ClearContents()

